$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated / reordered roster (18 players, rows 2..19).
$names = @(
    "Tyrese Maxey",
    "Stephen Curry",
    "Austin Reaves",
    "Darius Garland",
    "Kevin Durant",
    "Toumani Camara",
    "Karl-Anthony Towns",
    "Daniel Gafford",
    "Dereck Lively II",
    "Jarrett Allen",
    "Jalen Duren",
    "Mark Williams",
    "Tyrese Haliburton",
    "OG Anunoby",
    "Keegan Murray",
    "Franz Wagner",
    "Jalen Johnson",
    "Trey Murphy III"
)

$positions = @(
    "PG,SG",
    "PG,SG",
    "PG,SG",
    "PG",
    "SF,PF",
    "SF,PF",
    "PF,C",
    "PF,C",
    "C",
    "C",
    "C",
    "C",
    "PG,SG",
    "SF,PF",
    "SF,PF",
    "SF,PF",
    "SF,PF",
    "SF,PF"
)

$teams = @(
    "Philadelphia 76ers",
    "Golden State Warriors",
    "Los Angeles Lakers",
    "Cleveland Cavaliers",
    "Phoenix Suns",
    "Portland Trail Blazers",
    "New York Knicks",
    "Dallas Mavericks",
    "Dallas Mavericks",
    "Cleveland Cavaliers",
    "Detroit Pistons",
    "Charlotte Hornets",
    "Indiana Pacers",
    "New York Knicks",
    "Sacramento Kings",
    "Orlando Magic",
    "Atlanta Hawks",
    "New Orleans Pelicans"
)

# Fill column by column (A, then B, then C) to mirror the order in which
# the original edit populated the sheet.
for ($i = 0; $i -lt $names.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $names[$i]
}
for ($i = 0; $i -lt $positions.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $positions[$i]
}
for ($i = 0; $i -lt $teams.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $teams[$i]
}
